$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031341972238605
$ws.Range("D2").Value = 1.040362130796568
$ws.Range("E2").Value = 1.041184803876184
$ws.Range("F2").Value = 1.052746765887175
$ws.Range("I2").Value = 1.037328764366318
$ws.Range("J2").Value = 1.036478228737853
$ws.Range("K2").Value = 1.043144694831854
$ws.Range("L2").Value = 1.0439650361155
$ws.Range("M2").Value = 1.055494644700634
$ws.Range("N2").Value = 1.016154476248175
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032176259576357
$ws.Range("D3").Value = 1.040991524943093
$ws.Range("E3").Value = 1.041907733278212
$ws.Range("F3").Value = 1.053510512223387
$ws.Range("I3").Value = 1.037473554271788
$ws.Range("J3").Value = 1.036954937301612
$ws.Range("K3").Value = 1.043584973256931
$ws.Range("L3").Value = 1.044498771773278
$ws.Range("M3").Value = 1.056071424288792
$ws.Range("N3").Value = 1.016312732456598
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03271667834712
$ws.Range("D4").Value = 1.041398951770929
$ws.Range("E4").Value = 1.042376332613377
$ws.Range("F4").Value = 1.054005317416615
$ws.Range("I4").Value = 1.037565667005835
$ws.Range("J4").Value = 1.037263320288202
$ws.Range("K4").Value = 1.043869339311712
$ws.Range("L4").Value = 1.044844278450707
$ws.Range("M4").Value = 1.056444585346642
$ws.Range("N4").Value = 1.016415085955701
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032944007313189
$ws.Range("D5").Value = 1.041570271981262
$ws.Range("E5").Value = 1.042573525107426
$ws.Range("F5").Value = 1.054213477514688
$ws.Range("I5").Value = 1.037604012965903
$ws.Range("J5").Value = 1.037392944148297
$ws.Range("K5").Value = 1.043988759985448
$ws.Range("L5").Value = 1.044989562413284
$ws.Range("M5").Value = 1.056601447610802
$ws.Range("N5").Value = 1.016458103192827
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032982184825723
$ws.Range("D6").Value = 1.04159903956849
$ws.Range("E6").Value = 1.042606645875549
$ws.Range("F6").Value = 1.054248436893908
$ws.Range("I6").Value = 1.037610429214401
$ws.Range("J6").Value = 1.037414707332127
$ws.Range("K6").Value = 1.04400880376233
$ws.Range("L6").Value = 1.045013958118063
$ws.Range("M6").Value = 1.056627784565741
$ws.Range("N6").Value = 1.01646532524839
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032719715392062
$ws.Range("D7").Value = 1.041401240813216
$ws.Range("E7").Value = 1.042378966751912
$ws.Range("F7").Value = 1.05400809829888
$ws.Range("I7").Value = 1.037566180873801
$ws.Range("J7").Value = 1.037265052410225
$ws.Range("K7").Value = 1.043870935516768
$ws.Range("L7").Value = 1.044846219614553
$ws.Range("M7").Value = 1.056446681407631
$ws.Range("N7").Value = 1.016415660802568
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031623802616624
$ws.Range("D8").Value = 1.040574801806396
$ws.Range("E8").Value = 1.041428951197327
$ws.Range("F8").Value = 1.053004749922154
$ws.Range("I8").Value = 1.037378022558099
$ws.Range("J8").Value = 1.036639350240755
$ws.Range("K8").Value = 1.043293596496759
$ws.Range("L8").Value = 1.044145383845973
$ws.Range("M8").Value = 1.055689580179365
$ws.Range("N8").Value = 1.016207969501784
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029697175189799
$ws.Range("D9").Value = 1.039119873305124
$ws.Range("E9").Value = 1.039761235082623
$ws.Range("F9").Value = 1.051241488217078
$ws.Range("I9").Value = 1.037034434199645
$ws.Range("J9").Value = 1.035536235296458
$ws.Range("K9").Value = 1.04227230713103
$ws.Range("L9").Value = 1.042911591848709
$ws.Range("M9").Value = 1.054355134978896
$ws.Range("N9").Value = 1.01584163789356
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028415891733679
$ws.Range("D10").Value = 1.038150952545909
$ws.Range("E10").Value = 1.03865379749006
$ws.Range("F10").Value = 1.050069309791995
$ws.Range("I10").Value = 1.036797342541982
$ws.Range("J10").Value = 1.034800537515523
$ws.Range("K10").Value = 1.041588883582548
$ws.Range("L10").Value = 1.042089941077527
$ws.Range("M10").Value = 1.053465382195159
$ws.Range("N10").Value = 1.015597207935158
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027861846787001
$ws.Range("D11").Value = 1.037731668167951
$ws.Range("E11").Value = 1.038175326401658
$ws.Range("F11").Value = 1.049562559587439
$ws.Range("I11").Value = 1.036692785517918
$ws.Range("J11").Value = 1.034481920917943
$ws.Range("K11").Value = 1.041292363990074
$ws.Range("L11").Value = 1.041734384349065
$ws.Range("M11").Value = 1.053080102086177
$ws.Range("N11").Value = 1.015491323508937
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027656165606502
$ws.Range("D12").Value = 1.037575969014094
$ws.Range("E12").Value = 1.037997761667554
$ws.Range("F12").Value = 1.049374453937789
$ws.Range("I12").Value = 1.036653664620636
$ws.Range("J12").Value = 1.034363565682264
$ws.Range("K12").Value = 1.041182135846113
$ws.Range("L12").Value = 1.041602349859634
$ws.Range("M12").Value = 1.052936991954602
$ws.Range("N12").Value = 1.01545198709275
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027700279651045
$ws.Range("D13").Value = 1.03760936509117
$ws.Range("E13").Value = 1.03803584260065
$ws.Range("F13").Value = 1.049414797595558
$ws.Range("I13").Value = 1.036662069012279
$ws.Range("J13").Value = 1.034388953570047
$ws.Range("K13").Value = 1.04120578409595
$ws.Range("L13").Value = 1.041630670092368
$ws.Range("M13").Value = 1.052967689533055
$ws.Range("N13").Value = 1.015460425161571
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027844842744637
$ws.Range("D14").Value = 1.037718797164652
$ws.Range("E14").Value = 1.038160645555055
$ws.Range("F14").Value = 1.049547008175148
$ws.Range("I14").Value = 1.036689557555801
$ws.Range("J14").Value = 1.034472137769027
$ws.Range("K14").Value = 1.04128325427506
$ws.Range("L14").Value = 1.041723469611864
$ws.Range("M14").Value = 1.053068272546817
$ws.Range("N14").Value = 1.01548807207273
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027933928231228
$ws.Range("D15").Value = 1.037786227465466
$ws.Range("E15").Value = 1.038237562131835
$ws.Range("F15").Value = 1.049628483963448
$ws.Range("I15").Value = 1.036706456580358
$ws.Range("J15").Value = 1.034523389433015
$ws.Range("K15").Value = 1.041330974650251
$ws.Range("L15").Value = 1.041780651156976
$ws.Range("M15").Value = 1.053130245120891
$ws.Range("N15").Value = 1.01550510543655
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028452678014766
$ws.Range("D16").Value = 1.038178784813998
$ws.Range("E16").Value = 1.038685574472505
$ws.Range("F16").Value = 1.050102958398245
$ws.Range("I16").Value = 1.036804241803952
$ws.Range("J16").Value = 1.034821682015368
$ws.Range("K16").Value = 1.041608550250516
$ws.Range("L16").Value = 1.042113543046182
$ws.Range("M16").Value = 1.053490951878928
$ws.Range("N16").Value = 1.015604234236354
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.0287782805298
$ws.Range("D17").Value = 1.038425097924843
$ws.Range("E17").Value = 1.038966885081264
$ws.Range("F17").Value = 1.050400802250527
$ws.Range("I17").Value = 1.036865073246642
$ws.Range("J17").Value = 1.035008779412574
$ws.Range("K17").Value = 1.041782508458067
$ws.Range("L17").Value = 1.042322418249946
$ws.Range("M17").Value = 1.053717211915157
$ws.Range("N17").Value = 1.015666403501765
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028968272109268
$ws.Range("D18").Value = 1.038568793603928
$ws.Range("E18").Value = 1.039131070623717
$ws.Range("F18").Value = 1.050574607617422
$ws.Range("I18").Value = 1.036900372319437
$ws.Range("J18").Value = 1.035117904799522
$ws.Range("K18").Value = 1.041883918052792
$ws.Range("L18").Value = 1.04244427308197
$ws.Range("M18").Value = 1.053849184360571
$ws.Range("N18").Value = 1.015702661443032
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029033066681934
$ws.Range("D19").Value = 1.038617794367002
$ws.Range("E19").Value = 1.039187070861668
$ws.Range("F19").Value = 1.050633883936099
$ws.Range("I19").Value = 1.036912377344111
$ws.Range("J19").Value = 1.035155112785993
$ws.Range("K19").Value = 1.041918486342396
$ws.Range("L19").Value = 1.042485826023816
$ws.Range("M19").Value = 1.053894183296151
$ws.Range("N19").Value = 1.015715023713319
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028743338883694
$ws.Range("D20").Value = 1.038398668216796
$ws.Range("E20").Value = 1.038936692568872
$ws.Range("F20").Value = 1.050368838331678
$ws.Range("I20").Value = 1.036858565517989
$ws.Range("J20").Value = 1.034988706177062
$ws.Range("K20").Value = 1.041763850303246
$ws.Range("L20").Value = 1.042300005686604
$ws.Range("M20").Value = 1.053692936466641
$ws.Range("N20").Value = 1.015659733773324
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027802269306815
$ws.Range("D21").Value = 1.037686570978767
$ws.Range("E21").Value = 1.038123889742786
$ws.Range("F21").Value = 1.04950807201633
$ws.Range("I21").Value = 1.036681470690866
$ws.Range("J21").Value = 1.034447642273797
$ws.Range("K21").Value = 1.041260443643562
$ws.Range("L21").Value = 1.041696141484799
$ws.Range("M21").Value = 1.053038653340664
$ws.Range("N21").Value = 1.015479930914479
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027211252237808
$ws.Range("D22").Value = 1.037239089602529
$ws.Range("E22").Value = 1.037613779266255
$ws.Range("F22").Value = 1.048967592194381
$ws.Range("I22").Value = 1.036568482824345
$ws.Range("J22").Value = 1.034107414990901
$ws.Range("K22").Value = 1.040943426003033
$ws.Range("L22").Value = 1.04131667199097
$ws.Range("M22").Value = 1.052627280755243
$ws.Range("N22").Value = 1.015366845937587
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027524497419222
$ws.Range("D23").Value = 1.037476284313959
$ws.Range("E23").Value = 1.037884109582716
$ws.Range("F23").Value = 1.049254041995342
$ws.Range("I23").Value = 1.036628535106824
$ws.Range("J23").Value = 1.034287779148262
$ws.Range("K23").Value = 1.041111530512347
$ws.Range("L23").Value = 1.041517816110516
$ws.Range("M23").Value = 1.052845356450626
$ws.Range("N23").Value = 1.015426797654985
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02875912728373
$ws.Range("D24").Value = 1.038410610585946
$ws.Range("E24").Value = 1.038950334951234
$ws.Range("F24").Value = 1.050383281209727
$ws.Range("I24").Value = 1.036861506645612
$ws.Range("J24").Value = 1.034997776424872
$ws.Range("K24").Value = 1.041772281297457
$ws.Range("L24").Value = 1.042310132893124
$ws.Range("M24").Value = 1.053703905500841
$ws.Range("N24").Value = 1.015662747549978
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030194709490536
$ws.Range("D25").Value = 1.039495833802449
$ws.Range("E25").Value = 1.040191617308164
$ws.Range("F25").Value = 1.051696755922214
$ws.Range("I25").Value = 1.037124679878858
$ws.Range("J25").Value = 1.035821473562335
$ws.Range("K25").Value = 1.042536793250551
$ws.Range("L25").Value = 1.043230408646801
$ws.Range("M25").Value = 1.054700150215478
$ws.Range("N25").Value = 1.015936382291505
